# Natmi following Dr Hou advice
# Updates LR-pair statistics (Gas6-Mertk) for rows 2-17 to reflect
# recomputed ligand/receptor expressing-cell counts and derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  2 = @{ "E" = 3; "G" = 15.79785166666667; "H" = 47.39355500000001; "I" = 0.1445757693628457; "J" = 0.1445757693628457; "K" = 3; "M" = 7.731686; "N" = 23.195058; "O" = 0.1963057092861306; "P" = 0.1963057092861306; "Q" = 122.1440285612433; "R" = 1099.29625705119; "S" = 0.02838104895036145; "T" = 0.02838104895036146 }
  3 = @{ "E" = 3; "G" = 15.79785166666667; "H" = 47.39355500000001; "I" = 0.1445757693628457; "J" = 0.1445757693628457; "K" = 3; "M" = 2.796453333333333; "N" = 8.38936; "O" = 0.07100129972758389; "P" = 0.07100129972758387; "Q" = 44.17795495275556; "R" = 397.6015945748001; "S" = 0.01026506753387745; "T" = 0.01026506753387745 }
  4 = @{ "E" = 3; "G" = 15.79785166666667; "H" = 47.39355500000001; "I" = 0.1445757693628457; "J" = 0.1445757693628457; "K" = 3; "M" = 25.627053; "N" = 76.881159; "O" = 0.6506649152692259; "P" = 0.6506649152692259; "Q" = 404.852381947805; "R" = 3643.671437530245; "S" = 0.09407038072245916; "T" = 0.09407038072245917 }
  5 = @{ "E" = 3; "G" = 15.79785166666667; "H" = 47.39355500000001; "I" = 0.1445757693628457; "J" = 0.1445757693628457; "K" = 3; "M" = 3.230753333333334; "N" = 9.692260000000001; "O" = 0.08202807571705974; "P" = 0.08202807571705972; "Q" = 51.0389619315889; "R" = 459.3506573843001; "S" = 0.01185927215614767; "T" = 0.01185927215614767 }
  6 = @{ "E" = 3; "G" = 74.39645633333333; "H" = 223.189369; "I" = 0.6808473164079603; "J" = 0.6808473164079603; "K" = 3; "M" = 7.731686; "N" = 23.195058; "O" = 0.1963057092861306; "P" = 0.1963057092861306; "Q" = 575.2100398820446; "R" = 5176.890358938402; "S" = 0.1336542153630232; "T" = 0.1336542153630232 }
  7 = @{ "E" = 3; "G" = 74.39645633333333; "H" = 223.189369; "I" = 0.6808473164079603; "J" = 0.6808473164079603; "K" = 3; "M" = 2.796453333333333; "N" = 8.38936; "O" = 0.07100129972758389; "P" = 0.07100129972758387; "Q" = 208.0462183015378; "R" = 1872.41596471384; "S" = 0.04834104438100273; "T" = 0.04834104438100272 }
  8 = @{ "E" = 3; "G" = 74.39645633333333; "H" = 223.189369; "I" = 0.6808473164079603; "J" = 0.6808473164079603; "K" = 3; "M" = 25.627053; "N" = 76.881159; "O" = 0.6506649152692259; "P" = 0.6506649152692259; "Q" = 1906.561929466519; "R" = 17159.05736519867; "S" = 0.4430034614418653; "T" = 0.4430034614418653 }
  9 = @{ "E" = 3; "G" = 74.39645633333333; "H" = 223.189369; "I" = 0.6808473164079603; "J" = 0.6808473164079603; "K" = 3; "M" = 3.230753333333334; "N" = 9.692260000000001; "O" = 0.08202807571705974; "P" = 0.08202807571705972; "Q" = 240.3565992871045; "R" = 2163.20939358394; "S" = 0.05584859522206909; "T" = 0.05584859522206908 }
  10 = @{ "E" = 3; "G" = 16.36992; "H" = 49.10976; "I" = 0.1498111153557632; "J" = 0.1498111153557632; "K" = 3; "M" = 7.731686; "N" = 23.195058; "O" = 0.1963057092861306; "P" = 0.1963057092861306; "Q" = 126.56708128512; "R" = 1139.10373156608; "S" = 0.02940877725885941; "T" = 0.02940877725885941 }
  11 = @{ "E" = 3; "G" = 16.36992; "H" = 49.10976; "I" = 0.1498111153557632; "J" = 0.1498111153557632; "K" = 3; "M" = 2.796453333333333; "N" = 8.38936; "O" = 0.07100129972758389; "P" = 0.07100129972758387; "Q" = 45.7777173504; "R" = 411.9994561536; "S" = 0.01063678390389819; "T" = 0.01063678390389818 }
  12 = @{ "E" = 3; "G" = 16.36992; "H" = 49.10976; "I" = 0.1498111153557632; "J" = 0.1498111153557632; "K" = 3; "M" = 25.627053; "N" = 76.881159; "O" = 0.6506649152692259; "P" = 0.6506649152692259; "Q" = 419.51280744576; "R" = 3775.61526701184; "S" = 0.09747683667934587; "T" = 0.09747683667934587 }
  13 = @{ "E" = 3; "G" = 16.36992; "H" = 49.10976; "I" = 0.1498111153557632; "J" = 0.1498111153557632; "K" = 3; "M" = 3.230753333333334; "N" = 9.692260000000001; "O" = 0.08202807571705974; "P" = 0.08202807571705972; "Q" = 52.88717360640001; "R" = 475.9845624576001; "S" = 0.01228871751365971; "T" = 0.01228871751365971 }
  14 = @{ "E" = 3; "G" = 2.706168666666667; "H" = 8.118506; "I" = 0.02476579887343077; "J" = 0.02476579887343077; "K" = 3; "M" = 7.731686; "N" = 23.195058; "O" = 0.1963057092861306; "P" = 0.1963057092861306; "Q" = 20.92324639370533; "R" = 188.309217543348; "S" = 0.00486166771388648; "T" = 0.00486166771388648 }
  15 = @{ "E" = 3; "G" = 2.706168666666667; "H" = 8.118506; "I" = 0.02476579887343077; "J" = 0.02476579887343077; "K" = 3; "M" = 2.796453333333333; "N" = 8.38936; "O" = 0.07100129972758389; "P" = 0.07100129972758387; "Q" = 7.567674388462223; "R" = 68.10906949616; "S" = 0.001758403908805517; "T" = 0.001758403908805517 }
  16 = @{ "E" = 3; "G" = 2.706168666666667; "H" = 8.118506; "I" = 0.02476579887343077; "J" = 0.02476579887343077; "K" = 3; "M" = 25.627053; "N" = 76.881159; "O" = 0.6506649152692259; "P" = 0.6506649152692259; "Q" = 69.351127847606; "R" = 624.160150628454; "S" = 0.01611423642555552; "T" = 0.01611423642555552 }
  17 = @{ "E" = 3; "G" = 2.706168666666667; "H" = 8.118506; "I" = 0.02476579887343077; "J" = 0.02476579887343077; "K" = 3; "M" = 3.230753333333334; "N" = 9.692260000000001; "O" = 0.08202807571705974; "P" = 0.08202807571705972; "Q" = 8.742963440395556; "R" = 78.68667096356; "S" = 0.002031490825183252; "T" = 0.002031490825183252 }
}

foreach ($rowNum in $updates.Keys) {
  $rowValues = $updates[$rowNum]
  foreach ($col in $rowValues.Keys) {
    $ws.Range("$col$rowNum").Value = $rowValues[$col]
  }
}

Write-Output "Updated $($updates.Count) rows"
